# Update "want to go" counts (column F) for several rows across the
# 展览 / 演出 / 全部类型 sheets, per the latest site scrape (gh-pages
# output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 736
$ws1.Range("F10").Value = 1226
$ws1.Range("F11").Value = 649
$ws1.Range("F12").Value = 390
$ws1.Range("F13").Value = 517
$ws1.Range("F16").Value = 595
$ws1.Range("F18").Value = 363
$ws1.Range("F23").Value = 588

# --- 演出 (sheet 2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 329
$ws2.Range("F6").Value = 22

# --- 全部类型 (sheet 4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 329
$ws4.Range("F7").Value = 736
$ws4.Range("F14").Value = 1226
$ws4.Range("F15").Value = 649
$ws4.Range("F17").Value = 22
$ws4.Range("F18").Value = 390
$ws4.Range("F19").Value = 517
$ws4.Range("F23").Value = 595
$ws4.Range("F26").Value = 363
$ws4.Range("F37").Value = 588
